# Insert a new data row at row 194 (shifting existing rows 194-252 down to 195-253)
# and populate it with the new daily price record for Kiwi at
# "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("194:194").Insert()

$ws.Cells.Item(194, 1).Value  = 5
$ws.Cells.Item(194, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(194, 3).Value  = "Maule"
$ws.Cells.Item(194, 4).Value  = 44711
$ws.Cells.Item(194, 5).Value  = 7
$ws.Cells.Item(194, 6).Value  = "Fruta"
$ws.Cells.Item(194, 7).Value  = 100101
$ws.Cells.Item(194, 8).Value  = "Berries"
$ws.Cells.Item(194, 9).Value  = 100101007
$ws.Cells.Item(194, 10).Value = "Kiwi"
$ws.Cells.Item(194, 11).Value = "Hayward"
$ws.Cells.Item(194, 12).Value = "Primera"
$ws.Cells.Item(194, 13).Value = 200
$ws.Cells.Item(194, 14).Value = 7000
$ws.Cells.Item(194, 15).Value = 7000
$ws.Cells.Item(194, 16).Value = 7000
$ws.Cells.Item(194, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(194, 18).Value = "Provincia de Curic" + [char]0xF3
$ws.Cells.Item(194, 19).Value = 389
$ws.Cells.Item(194, 20).Value = 18
